$wb = $excel.ActiveWorkbook

# Rename sheets
$wsMarch = $wb.Worksheets.Item("March")
$wsMarch.Name = "mars"

$wsApril = $wb.Worksheets.Item("April")
$wsApril.Name = "april"

# Update Account value in April sheet (row 12, column E) from "ABC" to "gtdf"
$wsApril.Range("E12").Value = "gtdf"
